$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "58.042.59"
Set-TextCell $ws "E2" "  -4.16%  "
Set-TextCell $ws "D3" "2.605.41"
Set-TextCell $ws "E3" "  -3.52%  "
Set-TextCell $ws "E4" "  -0.12%  "
Set-TextCell $ws "D5" "516.76"
Set-TextCell $ws "E5" "  -1.51%  "
Set-TextCell $ws "D6" "142.14"
Set-TextCell $ws "E6" "  -2.17%  "
Set-TextCell $ws "E7" "  +0.31%  "
Set-TextCell $ws "E8" "  -1.66%  "
Set-TextCell $ws "E9" "  -0.68%  "
Set-TextCell $ws "E10" "  -2.79%  "
Set-TextCell $ws "D11" "0.337"
Set-TextCell $ws "E11" "  -0.44%  "
Set-TextCell $ws "E12" "  +0.93%  "
Set-TextCell $ws "D13" "3.063.30"
Set-TextCell $ws "E13" "  -3.65%  "
Set-TextCell $ws "D14" "58.027.10"
Set-TextCell $ws "E14" "  -4.18%  "
Set-TextCell $ws "D15" "20.90"
Set-TextCell $ws "E15" "  -1.63%  "
Set-TextCell $ws "E16" "  -1.82%  "
Set-TextCell $ws "D17" "2.622.07"
Set-TextCell $ws "E17" "  -3.64%  "
Set-TextCell $ws "E18" "  -2.36%  "
Set-TextCell $ws "D19" "333.64"
Set-TextCell $ws "E19" "  -3.37%  "
Set-TextCell $ws "D20" "10.34"
Set-TextCell $ws "E20" "  -2.57%  "
Set-TextCell $ws "E21" "  -2.94%  "
Set-TextCell $ws "E22" "  +0.24%  "
Set-TextCell $ws "D23" "63.92"
Set-TextCell $ws "E23" "  +0.89%  "
Set-TextCell $ws "E24" "  -1.44%  "
Set-TextCell $ws "E25" "  -2.67%  "
Set-TextCell $ws "D26" "1.00"
Set-TextCell $ws "E26" "  +0.70%  "
Set-TextCell $ws "E27" "  -2.33%  "
Set-TextCell $ws "E28" "  -4.28%  "
Set-TextCell $ws "E29" "  -3.16%  "
Set-TextCell $ws "E30" "  +0.15%  "
Set-TextCell $ws "D31" "1.57"
Set-TextCell $ws "E31" "  -1.47%  "
Set-TextCell $ws "B32" "Monero"
Set-TextCell $ws "C32" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D32" "150.79"
Set-TextCell $ws "E32" "  +0.79%  "
Set-TextCell $ws "B33" "EthereumClassic"
Set-TextCell $ws "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D33" "18.68"
Set-TextCell $ws "E33" "  -1.80%  "
Set-TextCell $ws "E34" "  -4.32%  "
Set-TextCell $ws "E35" "  -5.07%  "
Set-TextCell $ws "D36" "0.896"
Set-TextCell $ws "E36" "  -3.58%  "
Set-TextCell $ws "B37" "Fetch.AI"
Set-TextCell $ws "C37" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws "D37" "0.841"
Set-TextCell $ws "E37" "  -3.68%  "
Set-TextCell $ws "B38" "OKB"
Set-TextCell $ws "C38" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D38" "36.15"
Set-TextCell $ws "E38" "  -2.58%  "
Set-TextCell $ws "E39" "  -6.05%  "
Set-TextCell $ws "E40" "  -1.78%  "
Set-TextCell $ws "E41" "  +0.21%  "
Set-TextCell $ws "B42" "Stellar"
Set-TextCell $ws "C42" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D42" "0.0964"
Set-TextCell $ws "E42" "  -2.21%  "
Set-TextCell $ws "B43" "Mantle"
Set-TextCell $ws "C43" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws "D43" "0.595"
Set-TextCell $ws "E43" "  -2.82%  "
Set-TextCell $ws "D44" "267.88"
Set-TextCell $ws "E44" "  -5.03%  "
Set-TextCell $ws "E45" "  +1.20%  "
Set-TextCell $ws "B46" "Hedera"
Set-TextCell $ws "C46" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D46" "0.0531"
Set-TextCell $ws "E46" "  -1.24%  "
Set-TextCell $ws "B47" "EnergySwap"
Set-TextCell $ws "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D47" "19.04"
Set-TextCell $ws "E47" "  -5.11%  "
Set-TextCell $ws "D48" "2.030.88"
Set-TextCell $ws "E48" "  -5.20%  "
Set-TextCell $ws "E49" "  -1.36%  "
Set-TextCell $ws "D50" "4.62"
Set-TextCell $ws "E50" "  -6.07%  "
Set-TextCell $ws "D51" "18.10"
Set-TextCell $ws "E51" "  -4.98%  "
